# Atualizado por script em 05-11-2023 08:45
# Adds the two newly scraped ISL 2023-2024 matches as rows 32 and 33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing data row (31) down to
# the two new rows so the row-index column (A) and date column (E) keep
# the same cell styles (bold/bordered index column, date-formatted column).
$ws.Range("A31:V31").Copy()
$ws.Range("A32:V33").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Row 32
# ---------------------------------------------------------------------
$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = "india"
$ws.Cells.Item(32, 3).Value = "isl"
$ws.Cells.Item(32, 4).Value = "2023-2024"
$ws.Cells.Item(32, 5).Value = 45234.54166666666
$ws.Cells.Item(32, 6).Value = "Hyderabad"
$ws.Cells.Item(32, 7).Value = 1
$ws.Cells.Item(32, 8).Value = "Bengaluru FC"
$ws.Cells.Item(32, 9).Value = 1
$ws.Cells.Item(32, 10).Value = 1.95
$ws.Cells.Item(32, 11).Value = "31/10/2023 15:42"
$ws.Cells.Item(32, 12).Value = 2.18
$ws.Cells.Item(32, 13).Value = "04/11/2023 12:15"
$ws.Cells.Item(32, 14).Value = 3.48
$ws.Cells.Item(32, 15).Value = "31/10/2023 15:42"
$ws.Cells.Item(32, 16).Value = 3.78
$ws.Cells.Item(32, 17).Value = "04/11/2023 12:15"
$ws.Cells.Item(32, 18).Value = 3.9
$ws.Cells.Item(32, 19).Value = "31/10/2023 15:42"
$ws.Cells.Item(32, 20).Value = 3.11
$ws.Cells.Item(32, 21).Value = "04/11/2023 12:15"
$ws.Cells.Item(32, 22).Value = "https://www.betexplorer.com/football/india/isl/hyderabad-bengaluru-fc/O2e1xjG9/"

# ---------------------------------------------------------------------
# Row 33
# ---------------------------------------------------------------------
$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = "india"
$ws.Cells.Item(33, 3).Value = "isl"
$ws.Cells.Item(33, 4).Value = "2023-2024"
$ws.Cells.Item(33, 5).Value = 45234.64583333334
$ws.Cells.Item(33, 6).Value = "East Bengal"
$ws.Cells.Item(33, 7).Value = 1
$ws.Cells.Item(33, 8).Value = "Kerala Blasters"
$ws.Cells.Item(33, 9).Value = 2
$ws.Cells.Item(33, 10).Value = 2.7
$ws.Cells.Item(33, 11).Value = "28/10/2023 18:13"
$ws.Cells.Item(33, 12).Value = 2.32
$ws.Cells.Item(33, 13).Value = "04/11/2023 15:27"
$ws.Cells.Item(33, 14).Value = 3.35
$ws.Cells.Item(33, 15).Value = "28/10/2023 18:13"
$ws.Cells.Item(33, 16).Value = 3.29
$ws.Cells.Item(33, 17).Value = "04/11/2023 15:27"
$ws.Cells.Item(33, 18).Value = 2.6
$ws.Cells.Item(33, 19).Value = "28/10/2023 18:13"
$ws.Cells.Item(33, 20).Value = 3.23
$ws.Cells.Item(33, 21).Value = "04/11/2023 15:27"
$ws.Cells.Item(33, 22).Value = "https://www.betexplorer.com/football/india/isl/east-bengal-kerala-blasters/2gf5yAVF/"
